$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Row 3: "Test" header ---
$ws2.Range("A3").Value = "Test"

# --- Rows 4-6: Speed/Distance/Width labels ---
$ws2.Range("A4").Value = "Speed"
$ws2.Range("A5").Value = "Distance"
$ws2.Range("A6").Value = "Width"

# Column B, rows 4-10 are all text-formatted (range strings "10-20" etc.,
# two styled-but-empty cells, and the text value "46%").
$ws2.Range("B4:B10").NumberFormat = "@"

# Write B5 ("20-30") before B4 ("10-20") so the shared-string table
# picks up "20-30" first (matches the target workbook's string order).
$ws2.Range("B5").Value = "20-30"
$ws2.Range("B4").Value = "10-20"
$ws2.Range("B6").Value = "5-10"

# --- Row 9: "Config 1" header (column G) ---
$ws2.Range("G9").Value = "Config 1"

# --- "Avg error" label, first used in row 12 col A ---
$ws2.Range("A12").Value = "Avg error"

# --- Row 10: text "46%" in column B ---
$ws2.Range("B10").Value = "46%"

# --- Row 9: remaining config headers ---
$ws2.Range("A9").Value = "Config 3"
$ws2.Range("D9").Value = "Config 2"

# --- Numeric percentage cells (reuse existing 0% style) ---
$ws2.Range("A10").NumberFormat = "0%"
$ws2.Range("D10").NumberFormat = "0%"
$ws2.Range("G10").NumberFormat = "0%"
$ws2.Range("A11").NumberFormat = "0%"
$ws2.Range("D11").NumberFormat = "0%"
$ws2.Range("G11").NumberFormat = "0%"
$ws2.Range("B11").NumberFormat = "0%"

$ws2.Range("A10").Value = 0.1
$ws2.Range("D10").Value = 0.1
$ws2.Range("G10").Value = 0.1
$ws2.Range("A11").Value = 0.2
$ws2.Range("D11").Value = 0.2
$ws2.Range("G11").Value = 0.2
$ws2.Range("B11").Value = 0.03

# --- Row 12: remaining "Avg error" labels + plain numeric avg ---
$ws2.Range("D12").Value = "Avg error"
$ws2.Range("G12").Value = "Avg error"
$ws2.Range("B12").Value = 0.4

# --- Sheet2 print orientation ---
$ws2.PageSetup.Orientation = 1

# --- Sheet2 becomes the active sheet/tab with E14 selected ---
$ws2.Activate()
$ws2.Range("E14").Select() | Out-Null
